{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, and large numbers) inside specific bullet\n// paragraphs of the resume body, matching the target diff exactly.\n//\n// Strategy: locate each target paragraph by a unique leading/ containing\n// substring, then within that paragraph's scope, search for each metric\n// token (in left-to-right order) and apply bold + color formatting to the\n// matched sub-range. Word/Office.js automatically splits the run at the\n// match boundaries, which reproduces the same run structure as the diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map of: a substring that uniquely identifies the paragraph -> ordered\n// list of metric tokens to bold+color within that paragraph (in order of\n// appearance, matched with matchCase to avoid accidental partial hits).\nconst work = [\n  {\n    key: \"\u2022 Discovered systematic race coding errors\",\n    terms: [\"23%\", \"64%\"],\n  },\n  {\n    key: \"\u2022 Utilized advanced sampling methods\",\n    terms: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    key: \"\u2022 Trigonometric algorithm for boundary\",\n    terms: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    key: \"\u2022 Built real-time FEC analysis\",\n    terms: [\"$2\"],\n  },\n  {\n    key: \"\u2022 Modernized legacy ETL processes\",\n    terms: [\"57%\"],\n  },\n  {\n    key: \"\u2022 Algorithmic innovation: Pioneered\",\n    terms: [\"73.5%\"],\n  },\n  {\n    key: \"\u2022 $4.7M savings enabled\",\n    terms: [\"$4.7M\"],\n  },\n  {\n    key: \"\u2022 Platform impact: Built redistricting\",\n    terms: [\"12,847\"],\n  },\n];\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nfor (const item of work) {\n  // Find the first paragraph whose text starts with the target key.\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(item.key) === 0) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    continue;\n  }\n\n  for (const term of item.terms) {\n    const found = target.search(term, { matchCase: true });\n    found.load(\"text\");\n    await context.sync();\n    if (found.items.length === 0) {\n      continue;\n    }\n    const range = found.items[0];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, and large numbers) inside specific bullet\n# paragraphs of the resume body, matching the target diff exactly.\n#\n# Strategy: locate each target paragraph by a unique containing substring,\n# then within that paragraph's range, use Find.Execute to locate each\n# metric token (in left-to-right order) and apply Bold + Color formatting\n# to the found sub-range. Word splits the run at the match boundaries,\n# reproducing the same run structure as the diff.\n\nfunction HexToWordColor($hex) {\n    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)\n    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)\n    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$HIGHLIGHT_COLOR = HexToWordColor('2C3E50')\n\n$d = $word.ActiveDocument\n\n# Ordered list of (paragraph key substring, metric tokens to highlight).\n$work = @(\n    @{ Key = 'Discovered systematic race coding errors'; Terms = @('23%', '64%') },\n    @{ Key = 'Utilized advanced sampling methods'; Terms = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') },\n    @{ Key = 'Trigonometric algorithm for boundary'; Terms = @('73.5%', '$4.7M') },\n    @{ Key = 'Built real-time FEC analysis'; Terms = @('$2') },\n    @{ Key = 'Modernized legacy ETL processes'; Terms = @('57%') },\n    @{ Key = 'Algorithmic innovation: Pioneered'; Terms = @('73.5%') },\n    @{ Key = '$4.7M savings enabled'; Terms = @('$4.7M') },\n    @{ Key = 'Platform impact: Built redistricting'; Terms = @('12,847') }\n)\n\n$paraCount = $d.Paragraphs.Count\n\nforeach ($item in $work) {\n    $target = $null\n    for ($i = 1; $i -le $paraCount; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        if ($t.IndexOf($item.Key) -ge 0) {\n            $target = $p\n            break\n        }\n    }\n    if ($null -eq $target) {\n        continue\n    }\n\n    foreach ($term in $item.Terms) {\n        $r = $target.Range\n        $find = $r.Find\n        $find.ClearFormatting()\n        $found = $find.Execute($term)\n        if ($found) {\n            $r.Font.Bold = $true\n            $r.Font.Color = $HIGHLIGHT_COLOR\n        }\n    }\n}\n"}
